$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2566.625  # H70 was 2501
$ws.Cells.Item(70, 9).Value = 1700  # I70 was 1583.3334
$ws.Cells.Item(70, 10).Value = 2855.5  # J70 was 3877.5
$ws.Cells.Item(70, 11).Value = 5100  # K70 was 4750.0002
$ws.Cells.Item(70, 12).Value = 8566.5  # L70 was 11632.5
$ws.Cells.Item(70, 13).Value = -4830  # M70 was -4480.0002
$ws.Cells.Item(70, 14).Value = -9106.5  # N70 was -12172.5

$ws.Cells.Item(73, 8).Value = 2566.625  # H73 was 2501
$ws.Cells.Item(73, 9).Value = 1700  # I73 was 1583.3334
$ws.Cells.Item(73, 10).Value = 2855.5  # J73 was 3877.5
$ws.Cells.Item(73, 11).Value = 5100  # K73 was 4750.0002
$ws.Cells.Item(73, 12).Value = 8566.5  # L73 was 11632.5
$ws.Cells.Item(73, 13).Value = -4164  # M73 was -3814.0002
$ws.Cells.Item(73, 14).Value = -10438.5  # N73 was -13504.5

$ws.Cells.Item(93, 8).Value = 23026.629  # H93 was 23159.098
$ws.Cells.Item(93, 10).Value = 23026.629  # J93 was 23159.098
$ws.Cells.Item(93, 12).Value = 23026.629  # L93 was 23159.098
$ws.Cells.Item(93, 14).Value = -28018.629  # N93 was -28151.098

$ws.Cells.Item(98, 8).Value = 8121.5127  # H98 was 9052.709999999999
$ws.Cells.Item(98, 9).Value = 7219.706  # I98 was 8187.857
$ws.Cells.Item(98, 10).Value = 8818.362999999999  # J98 was 9764.941000000001
$ws.Cells.Item(98, 11).Value = 7219.706  # K98 was 8187.857
$ws.Cells.Item(98, 12).Value = 8818.362999999999  # L98 was 9764.941000000001
$ws.Cells.Item(98, 13).Value = -5721.706  # M98 was -6689.857
$ws.Cells.Item(98, 14).Value = -11814.363  # N98 was -12760.941

$ws.Cells.Item(112, 8).Value = 1281.1833  # H112 was 1292.8948
$ws.Cells.Item(112, 10).Value = 1294.4237  # J112 was 1307.0536
$ws.Cells.Item(112, 12).Value = 3883.2711  # L112 was 3921.1608
$ws.Cells.Item(112, 14).Value = -6099.2711  # N112 was -6137.1608

$ws.Cells.Item(122, 8).Value = 8121.5127  # H122 was 9052.709999999999
$ws.Cells.Item(122, 9).Value = 7219.706  # I122 was 8187.857
$ws.Cells.Item(122, 10).Value = 8818.362999999999  # J122 was 9764.941000000001
$ws.Cells.Item(122, 11).Value = 21659.118  # K122 was 24563.571
$ws.Cells.Item(122, 12).Value = 26455.089  # L122 was 29294.823
$ws.Cells.Item(122, 13).Value = -19209.118  # M122 was -22113.571
$ws.Cells.Item(122, 14).Value = -31355.089  # N122 was -34194.823

$ws.Cells.Item(132, 8).Value = 34833580  # H132 was 25902054
$ws.Cells.Item(132, 9).Value = 47625800  # I132 was 32262842
$ws.Cells.Item(132, 11).Value = 142877400  # K132 was 96788526
$ws.Cells.Item(132, 13).Value = -142874870  # M132 was -96785996

$ws.Cells.Item(138, 8).Value = 2759.8618  # H138 was 2738.4895
$ws.Cells.Item(138, 9).Value = 1537.2667  # I138 was 1466
$ws.Cells.Item(138, 10).Value = 2992  # J138 was 2992.9875
$ws.Cells.Item(138, 11).Value = 4611.800099999999  # K138 was 4398
$ws.Cells.Item(138, 12).Value = 8976  # L138 was 8978.962500000001
$ws.Cells.Item(138, 13).Value = 528.1999000000005  # M138 was 742
$ws.Cells.Item(138, 14).Value = -19256  # N138 was -19258.9625

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 25010000  # H6 was 16680000

$ws.Cells.Item(32, 8).Value = 5853.9307  # H32 was 6438.6094
$ws.Cells.Item(32, 9).Value = 4637.0444  # I32 was 5385.2705
$ws.Cells.Item(32, 11).Value = 4637.0444  # K32 was 5385.2705
$ws.Cells.Item(32, 13).Value = -4350.0444  # M32 was -5098.2705

$ws.Cells.Item(45, 8).Value = 3496  # H45 was 3198.4
$ws.Cells.Item(45, 9).Value = 5740  # I45 was 3997.3333
$ws.Cells.Item(45, 11).Value = 5740  # K45 was 3997.3333
$ws.Cells.Item(45, 13).Value = -5363  # M45 was -3620.3333

$ws.Cells.Item(61, 8).Value = 1719.2858  # H61 was 1774.75
$ws.Cells.Item(61, 9).Value = 1705.8334  # I61 was 1774.75
$ws.Cells.Item(61, 10).Value = 1800  # J61 was 0
$ws.Cells.Item(61, 11).Value = 1705.8334  # K61 was 1774.75
$ws.Cells.Item(61, 12).Value = 1800  # L61 was 0
$ws.Cells.Item(61, 13).Value = -1493.8334  # M61 was -1562.75
$ws.Cells.Item(61, 14).Value = -2224  # N61 was None

$ws.Cells.Item(74, 8).Value = 3018.8696  # H74 was 2524.4138
$ws.Cells.Item(74, 9).Value = 2648.2354  # I74 was 2108.4348
$ws.Cells.Item(74, 10).Value = 4069  # J74 was 4119
$ws.Cells.Item(74, 11).Value = 2648.2354  # K74 was 2108.4348
$ws.Cells.Item(74, 12).Value = 4069  # L74 was 4119
$ws.Cells.Item(74, 13).Value = -1774.2354  # M74 was -1234.4348
$ws.Cells.Item(74, 14).Value = -5817  # N74 was -5867

$ws.Cells.Item(77, 8).Value = 3018.8696  # H77 was 2524.4138
$ws.Cells.Item(77, 9).Value = 2648.2354  # I77 was 2108.4348
$ws.Cells.Item(77, 10).Value = 4069  # J77 was 4119
$ws.Cells.Item(77, 11).Value = 13241.177  # K77 was 10542.174
$ws.Cells.Item(77, 12).Value = 20345  # L77 was 20595
$ws.Cells.Item(77, 13).Value = -8873.177  # M77 was -6174.173999999999
$ws.Cells.Item(77, 14).Value = -29081  # N77 was -29331

$ws.Cells.Item(106, 8).Value = 35000  # H106 was 50000
$ws.Cells.Item(106, 10).Value = 35000  # J106 was 50000
$ws.Cells.Item(106, 12).Value = 35000  # L106 was 50000
$ws.Cells.Item(106, 14).Value = -37524  # N106 was -52524

$ws.Cells.Item(132, 8).Value = 3170.5833  # H132 was 2315.3157
$ws.Cells.Item(132, 9).Value = 1247  # I132 was 1035.2142
$ws.Cells.Item(132, 10).Value = 4544.5713  # J132 was 5899.6
$ws.Cells.Item(132, 11).Value = 3741  # K132 was 3105.6426
$ws.Cells.Item(132, 12).Value = 13633.7139  # L132 was 17698.8
$ws.Cells.Item(132, 13).Value = -1211  # M132 was -575.6425999999997
$ws.Cells.Item(132, 14).Value = -18693.7139  # N132 was -22758.8

$ws.Cells.Item(136, 8).Value = 1719.2858  # H136 was 1774.75
$ws.Cells.Item(136, 9).Value = 1705.8334  # I136 was 1774.75
$ws.Cells.Item(136, 10).Value = 1800  # J136 was 0
$ws.Cells.Item(136, 11).Value = 5117.5002  # K136 was 5324.25
$ws.Cells.Item(136, 12).Value = 5400  # L136 was 0
$ws.Cells.Item(136, 13).Value = -2567.5002  # M136 was -2774.25
$ws.Cells.Item(136, 14).Value = -10500  # N136 was None

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(95, 8).Value = 33888.89  # H95 was 33000
$ws.Cells.Item(95, 10).Value = 33888.89  # J95 was 33000
$ws.Cells.Item(95, 12).Value = 33888.89  # L95 was 33000
$ws.Cells.Item(95, 14).Value = -39380.89  # N95 was -38492

$ws.Cells.Item(134, 8).Value = 4377.9165  # H134 was 2786.1904
$ws.Cells.Item(134, 9).Value = 1815.125  # I134 was 1527.5555
$ws.Cells.Item(134, 10).Value = 9503.5  # J134 was 10338
$ws.Cells.Item(134, 11).Value = 5445.375  # K134 was 4582.666499999999
$ws.Cells.Item(134, 12).Value = 28510.5  # L134 was 31014
$ws.Cells.Item(134, 13).Value = -2910.375  # M134 was -2047.666499999999
$ws.Cells.Item(134, 14).Value = -33580.5  # N134 was -36084

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3280.9546  # H31 was 2940.4033
$ws.Cells.Item(31, 9).Value = 2136.4285  # I31 was 1597.341
$ws.Cells.Item(31, 10).Value = 5283.875  # J31 was 6223.4443
$ws.Cells.Item(31, 11).Value = 2136.4285  # K31 was 1597.341
$ws.Cells.Item(31, 12).Value = 5283.875  # L31 was 6223.4443
$ws.Cells.Item(31, 13).Value = -1841.4285  # M31 was -1302.341
$ws.Cells.Item(31, 14).Value = -5873.875  # N31 was -6813.4443

$ws.Cells.Item(34, 8).Value = 3280.9546  # H34 was 2940.4033
$ws.Cells.Item(34, 9).Value = 2136.4285  # I34 was 1597.341
$ws.Cells.Item(34, 10).Value = 5283.875  # J34 was 6223.4443
$ws.Cells.Item(34, 11).Value = 2136.4285  # K34 was 1597.341
$ws.Cells.Item(34, 12).Value = 5283.875  # L34 was 6223.4443
$ws.Cells.Item(34, 13).Value = -1934.4285  # M34 was -1395.341
$ws.Cells.Item(34, 14).Value = -5687.875  # N34 was -6627.4443

$ws.Cells.Item(62, 8).Value = 3225  # H62 was 3266.8333
$ws.Cells.Item(62, 9).Value = 2966.6667  # I62 was 3001.6667
$ws.Cells.Item(62, 10).Value = 4000  # J62 was 3532
$ws.Cells.Item(62, 11).Value = 2966.6667  # K62 was 3001.6667
$ws.Cells.Item(62, 12).Value = 4000  # L62 was 3532
$ws.Cells.Item(62, 13).Value = -2342.6667  # M62 was -2377.6667
$ws.Cells.Item(62, 14).Value = -5248  # N62 was -4780

$ws.Cells.Item(65, 8).Value = 3225  # H65 was 3266.8333
$ws.Cells.Item(65, 9).Value = 2966.6667  # I65 was 3001.6667
$ws.Cells.Item(65, 10).Value = 4000  # J65 was 3532
$ws.Cells.Item(65, 11).Value = 14833.3335  # K65 was 15008.3335
$ws.Cells.Item(65, 12).Value = 20000  # L65 was 17660
$ws.Cells.Item(65, 13).Value = -11713.3335  # M65 was -11888.3335
$ws.Cells.Item(65, 14).Value = -26240  # N65 was -23900

$ws.Cells.Item(103, 8).Value = 14904.8  # H103 was 19874.77
$ws.Cells.Item(103, 9).Value = 6256  # I103 was 7131.125
$ws.Cells.Item(103, 10).Value = 49500  # J103 was 40264.6
$ws.Cells.Item(103, 11).Value = 6256  # K103 was 7131.125
$ws.Cells.Item(103, 12).Value = 49500  # L103 was 40264.6
$ws.Cells.Item(103, 13).Value = -5084  # M103 was -5959.125
$ws.Cells.Item(103, 14).Value = -51844  # N103 was -42608.6

$ws.Cells.Item(132, 8).Value = 3927.625  # H132 was 4404.5386
$ws.Cells.Item(132, 9).Value = 2241.7273  # I132 was 2384.5
$ws.Cells.Item(132, 11).Value = 6725.1819  # K132 was 7153.5
$ws.Cells.Item(132, 13).Value = -4195.1819  # M132 was -4623.5

$ws.Cells.Item(134, 8).Value = 7992.0527  # H134 was 7495.7617
$ws.Cells.Item(134, 9).Value = 8523.267  # I134 was 9614.846
$ws.Cells.Item(134, 10).Value = 6000  # J134 was 4052.25
$ws.Cells.Item(134, 11).Value = 25569.801  # K134 was 28844.538
$ws.Cells.Item(134, 12).Value = 18000  # L134 was 12156.75
$ws.Cells.Item(134, 13).Value = -23034.801  # M134 was -26309.538
$ws.Cells.Item(134, 14).Value = -23070  # N134 was -17226.75

$ws.Cells.Item(137, 8).Value = 50446.668  # H137 was 50780
$ws.Cells.Item(137, 10).Value = 50446.668  # J137 was 50780
$ws.Cells.Item(137, 12).Value = 50446.668  # L137 was 50780
$ws.Cells.Item(137, 14).Value = -60646.668  # N137 was -60980

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 535096.5  # H5 was 581554.6
$ws.Cells.Item(5, 9).Value = 399.6  # I5 was 409.92856
$ws.Cells.Item(5, 10).Value = 1337141.8  # J5 was 1485557.5
$ws.Cells.Item(5, 11).Value = 1198.8  # K5 was 1229.78568
$ws.Cells.Item(5, 12).Value = 4011425.4  # L5 was 4456672.5
$ws.Cells.Item(5, 13).Value = -1086.8  # M5 was -1117.78568
$ws.Cells.Item(5, 14).Value = -4011649.4  # N5 was -4456896.5

$ws.Cells.Item(14, 8).Value = 106.38461  # H14 was 116.90909
$ws.Cells.Item(14, 9).Value = 106.38461  # I14 was 116.90909
$ws.Cells.Item(14, 11).Value = 319.15383  # K14 was 350.72727
$ws.Cells.Item(14, 13).Value = -146.15383  # M14 was -177.72727

$ws.Cells.Item(133, 8).Value = 4960.6665  # H133 was 5542.8
$ws.Cells.Item(133, 9).Value = 4992.8  # I133 was 5728.5
$ws.Cells.Item(133, 11).Value = 14978.4  # K133 was 17185.5
$ws.Cells.Item(133, 13).Value = -9918.400000000001  # M133 was -12125.5

$ws.Cells.Item(135, 8).Value = 535096.5  # H135 was 581554.6
$ws.Cells.Item(135, 9).Value = 399.6  # I135 was 409.92856
$ws.Cells.Item(135, 10).Value = 1337141.8  # J135 was 1485557.5
$ws.Cells.Item(135, 11).Value = 3596.4  # K135 was 3689.35704
$ws.Cells.Item(135, 12).Value = 12034276.2  # L135 was 13370017.5
$ws.Cells.Item(135, 13).Value = -1061.4  # M135 was -1154.35704
$ws.Cells.Item(135, 14).Value = -12039346.2  # N135 was -13375087.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(82, 8).Value = 37859.6  # H82 was 30000
$ws.Cells.Item(82, 9).Value = 19649  # I82 was 30000
$ws.Cells.Item(82, 10).Value = 50000  # J82 was 0
$ws.Cells.Item(82, 11).Value = 19649  # K82 was 30000
$ws.Cells.Item(82, 12).Value = 50000  # L82 was 0
$ws.Cells.Item(82, 13).Value = -19266  # M82 was -29617
$ws.Cells.Item(82, 14).Value = -50766  # N82 was None

$ws.Cells.Item(85, 8).Value = 37859.6  # H85 was 30000
$ws.Cells.Item(85, 9).Value = 19649  # I85 was 30000
$ws.Cells.Item(85, 10).Value = 50000  # J85 was 0
$ws.Cells.Item(85, 11).Value = 19649  # K85 was 30000
$ws.Cells.Item(85, 12).Value = 50000  # L85 was 0
$ws.Cells.Item(85, 13).Value = -18323  # M85 was -28674
$ws.Cells.Item(85, 14).Value = -52652  # N85 was None

$ws.Cells.Item(126, 8).Value = 3354.59  # H126 was 3388.59
$ws.Cells.Item(126, 9).Value = 2837.808  # I126 was 2855
$ws.Cells.Item(126, 10).Value = 4751.815  # J126 was 4760.6787
$ws.Cells.Item(126, 11).Value = 8513.423999999999  # K126 was 8565
$ws.Cells.Item(126, 12).Value = 14255.445  # L126 was 14282.0361
$ws.Cells.Item(126, 13).Value = -6043.423999999999  # M126 was -6095
$ws.Cells.Item(126, 14).Value = -19195.445  # N126 was -19222.0361

$ws.Cells.Item(132, 8).Value = 8999.666999999999  # H132 was 6287.143
$ws.Cells.Item(132, 9).Value = 0  # I132 was 3670.6667
$ws.Cells.Item(132, 10).Value = 8999.666999999999  # J132 was 8249.5
$ws.Cells.Item(132, 11).Value = 0  # K132 was 11012.0001
$ws.Cells.Item(132, 12).ClearContents()  # L132 was 24748.5, removed
$ws.Cells.Item(132, 13).Value = 26999.001  # M132 was -8482.000100000001
$ws.Cells.Item(132, 14).Value = -32059.001  # N132 was -29808.5

$ws.Cells.Item(141, 8).Value = 37955  # H141 was 40330
$ws.Cells.Item(141, 10).Value = 38606.668  # J141 was 41773.332
$ws.Cells.Item(141, 12).Value = 38606.668  # L141 was 41773.332
$ws.Cells.Item(141, 14).Value = -48966.668  # N141 was -52133.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 4727.5  # H9 was 5352.143
$ws.Cells.Item(9, 9).Value = 205  # I9 was 185.8
$ws.Cells.Item(9, 10).Value = 9250  # J9 was 8222.333000000001
$ws.Cells.Item(9, 11).Value = 205  # K9 was 185.8
$ws.Cells.Item(9, 12).Value = 9250  # L9 was 8222.333000000001
$ws.Cells.Item(9, 13).Value = 19  # M9 was 38.19999999999999
$ws.Cells.Item(9, 14).Value = -9698  # N9 was -8670.333000000001

$ws.Cells.Item(26, 8).Value = 16002.667  # H26 was 16336
$ws.Cells.Item(26, 10).Value = 21499.5  # J26 was 21999.5
$ws.Cells.Item(26, 12).Value = 21499.5  # L26 was 21999.5
$ws.Cells.Item(26, 14).Value = -22089.5  # N26 was -22589.5

$ws.Cells.Item(132, 8).Value = 5883.654  # H132 was 5984.185
$ws.Cells.Item(132, 9).Value = 1943.1111  # I132 was 2076.389
$ws.Cells.Item(132, 10).Value = 14749.875  # J132 was 13799.777
$ws.Cells.Item(132, 11).Value = 5829.3333  # K132 was 6229.167
$ws.Cells.Item(132, 12).Value = 44249.625  # L132 was 41399.331
$ws.Cells.Item(132, 13).Value = -3299.3333  # M132 was -3699.167
$ws.Cells.Item(132, 14).Value = -49309.625  # N132 was -46459.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 5242.45  # H113 was 4551.7827
$ws.Cells.Item(113, 9).Value = 6121  # I113 was 7899.846
$ws.Cells.Item(113, 10).Value = 264  # J113 was 199.3
$ws.Cells.Item(113, 11).Value = 18363  # K113 was 23699.538
$ws.Cells.Item(113, 12).Value = 792  # L113 was 597.9000000000001
$ws.Cells.Item(113, 13).Value = -16193  # M113 was -21529.538
$ws.Cells.Item(113, 14).Value = -5132  # N113 was -4937.9

$ws.Cells.Item(132, 8).Value = 13895047  # H132 was 15880025
$ws.Cells.Item(132, 9).Value = 7539  # I132 was 9339.416999999999
$ws.Cells.Item(132, 10).Value = 37040896  # J132 was 37040940
$ws.Cells.Item(132, 11).Value = 22617  # K132 was 28018.251
$ws.Cells.Item(132, 12).Value = 111122688  # L132 was 111122820
$ws.Cells.Item(132, 13).Value = -20087  # M132 was -25488.251
$ws.Cells.Item(132, 14).Value = -111127748  # N132 was -111127880

$ws.Cells.Item(136, 8).Value = 6896.5938  # H136 was 7999.815
$ws.Cells.Item(136, 9).Value = 7368.1113  # I136 was 9840.77
$ws.Cells.Item(136, 11).Value = 22104.3339  # K136 was 29522.31
$ws.Cells.Item(136, 13).Value = -19554.3339  # M136 was -26972.31
